# Fruta / hortaliza, semanal
# Insert 4 new weekly rows (date 44448) right after row 458 for
# "Terminal La Palmera de La Serena - Naranja", shifting the existing
# rows 459-472 down to 463-476.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 459 (shifts 459:472 -> 463:476)
$ws.Rows.Item(459).Insert()
$ws.Rows.Item(459).Insert()
$ws.Rows.Item(459).Insert()
$ws.Rows.Item(459).Insert()

# Common columns shared by every data row in this block
$marketId = 8
$market = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad = "$/bins (400 kilos)"
$origen = "Provincia de Limarí"
$kgUnidad = 400

function Set-DataRow {
    param($row, $fecha, $variedad, $calidad, $volumen, $pmin, $pmax, $pprom, $precioKg)

    $ws.Cells.Item($row, 1).Value = $marketId
    $ws.Cells.Item($row, 2).Value = $market
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pprom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DataRow 459 44448 "Lane Late"   "Primera" 20 125000 130000 127500 319
Set-DataRow 460 44448 "Lane Late"   "Segunda" 20  90000 100000  95000 238
Set-DataRow 461 44448 "Navel Late"  "Primera" 20 125000 130000 127500 319
Set-DataRow 462 44448 "Navel Late"  "Segunda" 20  90000 100000  95000 238
